$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, pushing existing rows 36-39 down to 37-40
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly data record
$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = 44858
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = 100112028
$ws.Range("G36").Value = "Sandia"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Segunda"
$ws.Range("J36").Value = 700
$ws.Range("K36").Value = 730
$ws.Range("L36").Value = 750
$ws.Range("M36").Value = 740
$ws.Range("N36").Value = "$/kilo (volumen en unidades)"
$ws.Range("O36").Value = "Perú"
$ws.Range("P36").Value = 740
$ws.Range("Q36").Value = 1
$ws.Range("R36").Value = "Hortaliza"
